$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33518.332
$ws.Range("J3").Value = 33518.332
$ws.Range("L3").Value = 33518.332
$ws.Range("N3").Value = -33746.332
$ws.Range("H100").Value = 2270.3225
$ws.Range("I100").Value = 1810.4348
$ws.Range("J100").Value = 3592.5
$ws.Range("K100").Value = 1810.4348
$ws.Range("L100").Value = 3592.5
$ws.Range("M100").Value = -1269.4348
$ws.Range("N100").Value = -4674.5
$ws.Range("H102").Value = 33518.332
$ws.Range("J102").Value = 33518.332
$ws.Range("L102").Value = 33518.332
$ws.Range("N102").Value = -40008.332
$ws.Range("H116").Value = 5082.1904
$ws.Range("I116").Value = 4860
$ws.Range("J116").Value = 5284.1816
$ws.Range("K116").Value = 4860
$ws.Range("L116").Value = 5284.1816
$ws.Range("M116").Value = -1418
$ws.Range("N116").Value = -12168.1816
$ws.Range("H131").Value = 3406.7222
$ws.Range("I131").Value = 2773.0833
$ws.Range("K131").Value = 8319.249899999999
$ws.Range("M131").Value = -3279.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1329.6154
$ws.Range("I45").Value = 1024.5667
$ws.Range("J45").Value = 2346.4443
$ws.Range("K45").Value = 1024.5667
$ws.Range("L45").Value = 2346.4443
$ws.Range("M45").Value = -647.5667000000001
$ws.Range("N45").Value = -3100.4443
$ws.Range("H97").Value = 497.53845
$ws.Range("I97").Value = 384.36365
$ws.Range("J97").Value = 1120
$ws.Range("K97").Value = 384.36365
$ws.Range("L97").Value = 1120
$ws.Range("M97").Value = 111.63635
$ws.Range("N97").Value = -2112
$ws.Range("H102").Value = 3673.625
$ws.Range("I102").Value = 2881.5
$ws.Range("J102").Value = 6050
$ws.Range("K102").Value = 2881.5
$ws.Range("L102").Value = 6050
$ws.Range("M102").Value = -1259.5
$ws.Range("N102").Value = -9294
$ws.Range("H132").Value = 1823.0667
$ws.Range("I132").Value = 1519.5532
$ws.Range("J132").Value = 2920.3845
$ws.Range("K132").Value = 4558.6596
$ws.Range("L132").Value = 8761.1535
$ws.Range("M132").Value = -2028.6596
$ws.Range("N132").Value = -13821.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 24397.541
$ws.Range("I86").Value = 2430.0557
$ws.Range("J86").Value = 90300
$ws.Range("K86").Value = 2430.0557
$ws.Range("L86").Value = 90300
$ws.Range("M86").Value = -1307.0557
$ws.Range("N86").Value = -92546
$ws.Range("H89").Value = 24397.541
$ws.Range("I89").Value = 2430.0557
$ws.Range("J89").Value = 90300
$ws.Range("K89").Value = 12150.2785
$ws.Range("L89").Value = 451500
$ws.Range("M89").Value = -6534.2785
$ws.Range("N89").Value = -462732
$ws.Range("H99").Value = 2166.5
$ws.Range("I99").Value = 1544.2142
$ws.Range("J99").Value = 3618.5
$ws.Range("K99").Value = 1544.2142
$ws.Range("L99").Value = 3618.5
$ws.Range("M99").Value = -46.21419999999989
$ws.Range("N99").Value = -6614.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10872850
$ws.Range("I58").Value = 1848.3793
$ws.Range("K58").Value = 1848.3793
$ws.Range("M58").Value = -1645.3793
$ws.Range("H99").Value = 2467.3845
$ws.Range("I99").Value = 1841.7142
$ws.Range("J99").Value = 3197.3333
$ws.Range("K99").Value = 1841.7142
$ws.Range("L99").Value = 3197.3333
$ws.Range("M99").Value = -343.7141999999999
$ws.Range("N99").Value = -6193.3333
$ws.Range("H107").Value = 2286.2666
$ws.Range("I107").Value = 642.3333
$ws.Range("J107").Value = 4752.1665
$ws.Range("K107").Value = 642.3333
$ws.Range("L107").Value = 4752.1665
$ws.Range("M107").Value = 1277.6667
$ws.Range("N107").Value = -8592.166499999999
$ws.Range("H126").Value = 2467.3845
$ws.Range("I126").Value = 1841.7142
$ws.Range("J126").Value = 3197.3333
$ws.Range("K126").Value = 5525.142599999999
$ws.Range("L126").Value = 9591.999899999999
$ws.Range("M126").Value = -3055.142599999999
$ws.Range("N126").Value = -14531.9999
$ws.Range("H132").Value = 1938.717
$ws.Range("I132").Value = 1631.1875
$ws.Range("J132").Value = 2407.3333
$ws.Range("K132").Value = 4893.5625
$ws.Range("L132").Value = 7221.999899999999
$ws.Range("M132").Value = -2363.5625
$ws.Range("N132").Value = -12281.9999
$ws.Range("H134").Value = 2791.3
$ws.Range("I134").Value = 1559
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 4677
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -2142
$ws.Range("N134").Value = -22069.9995
$ws.Range("H136").Value = 10872850
$ws.Range("I136").Value = 1848.3793
$ws.Range("K136").Value = 5545.1379
$ws.Range("M136").Value = -2995.1379

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 14337.333
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61498
$ws.Range("H66").Value = 14337.333
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -187488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 37923.332
$ws.Range("I48").Value = 15980
$ws.Range("J48").Value = 42312
$ws.Range("K48").Value = 15980
$ws.Range("L48").Value = 42312
$ws.Range("M48").Value = -15495
$ws.Range("N48").Value = -43282
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H97").Value = 768.3333
$ws.Range("I97").Value = 387.95456
$ws.Range("J97").Value = 2442
$ws.Range("K97").Value = 387.95456
$ws.Range("L97").Value = 2442
$ws.Range("M97").Value = 108.04544
$ws.Range("N97").Value = -3434
$ws.Range("H102").Value = 29998.277
$ws.Range("I102").Value = 1885.96
$ws.Range("J102").Value = 93889.91
$ws.Range("K102").Value = 1885.96
$ws.Range("L102").Value = 93889.91
$ws.Range("M102").Value = -263.96
$ws.Range("N102").Value = -97133.91
$ws.Range("H122").Value = 3333.3235
$ws.Range("I122").Value = 2158.7144
$ws.Range("J122").Value = 5230.769
$ws.Range("K122").Value = 6476.1432
$ws.Range("L122").Value = 15692.307
$ws.Range("M122").Value = -4026.1432
$ws.Range("N122").Value = -20592.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1680.5769
$ws.Range("I7").Value = 1107.3636
$ws.Range("J7").Value = 2100.9333
$ws.Range("K7").Value = 1107.3636
$ws.Range("L7").Value = 2100.9333
$ws.Range("M7").Value = -995.3635999999999
$ws.Range("N7").Value = -2324.9333
$ws.Range("H93").Value = 1598.8462
$ws.Range("I93").Value = 922.625
$ws.Range("J93").Value = 2680.8
$ws.Range("K93").Value = 922.625
$ws.Range("L93").Value = 2680.8
$ws.Range("M93").Value = 325.375
$ws.Range("N93").Value = -5176.8
$ws.Range("H126").Value = 1680.5769
$ws.Range("I126").Value = 1107.3636
$ws.Range("J126").Value = 2100.9333
$ws.Range("K126").Value = 3322.0908
$ws.Range("L126").Value = 6302.7999
$ws.Range("M126").Value = -852.0907999999999
$ws.Range("N126").Value = -11242.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 771814.9
$ws.Range("I122").Value = 911053.8
$ws.Range("J122").Value = 6000.5
$ws.Range("K122").Value = 2733161.4
$ws.Range("L122").Value = 18001.5
$ws.Range("M122").Value = -2730711.4
$ws.Range("N122").Value = -22901.5
$ws.Range("H126").Value = 3847964.2
$ws.Range("I126").Value = 1426.7646
$ws.Range("J126").Value = 11113646
$ws.Range("K126").Value = 4280.293799999999
$ws.Range("L126").Value = 33340938
$ws.Range("M126").Value = -1810.293799999999
$ws.Range("N126").Value = -33345878
